$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.051203666666667
$ws.Range("H2").Value = 18.153611
$ws.Range("I2").Value = 0.07608037240065801
$ws.Range("J2").Value = 0.07775008964215516
$ws.Range("M2").Value = 98.170451
$ws.Range("N2").Value = 294.511353
$ws.Range("O2").Value = 0.4110278868558457
$ws.Range("P2").Value = 0.4243928215400083
$ws.Range("Q2").Value = 594.0493930495204
$ws.Range("R2").Value = 5346.444537445683
$ws.Range("S2").Value = 0.03127115469904827
$ws.Range("T2").Value = 0.0329965799182228

$ws.Range("G3").Value = 6.051203666666667
$ws.Range("H3").Value = 18.153611
$ws.Range("I3").Value = 0.07608037240065801
$ws.Range("J3").Value = 0.07775008964215516
$ws.Range("O3").Value = 0.1453748421312515
$ws.Range("P3").Value = 0.150101833491052
$ws.Range("Q3").Value = 210.1070012386473
$ws.Range("R3").Value = 1890.963011147826
$ws.Range("S3").Value = 0.01106017212703248
$ws.Range("T3").Value = 0.01167043100938114

$ws.Range("G4").Value = 6.051203666666667
$ws.Range("H4").Value = 18.153611
$ws.Range("I4").Value = 0.07608037240065801
$ws.Range("J4").Value = 0.07775008964215516
$ws.Range("M4").Value = 38.82199566666667
$ws.Range("N4").Value = 116.465987
$ws.Range("O4").Value = 0.1625430328561575
$ws.Range("P4").Value = 0.1678282630971171
$ws.Range("Q4").Value = 234.9198025254508
$ws.Range("R4").Value = 2114.278222729057
$ws.Range("S4").Value = 0.01236633447082885
$ws.Range("T4").Value = 0.01304866250028806

$ws.Range("G5").Value = 6.051203666666667
$ws.Range("H5").Value = 18.153611
$ws.Range("I5").Value = 0.07608037240065801
$ws.Range("J5").Value = 0.07775008964215516
$ws.Range("M5").Value = 22.5647
$ws.Range("N5").Value = 45.1294
$ws.Range("O5").Value = 0.09447568860141126
$ws.Range("P5").Value = 0.06503176602637677
$ws.Range("Q5").Value = 136.5435953772333
$ws.Range("R5").Value = 819.2615722634001
$ws.Range("S5").Value = 0.00718774557160397
$ws.Range("T5").Value = 0.005056225638138454

$ws.Range("G6").Value = 6.051203666666667
$ws.Range("H6").Value = 18.153611
$ws.Range("I6").Value = 0.07608037240065801
$ws.Range("J6").Value = 0.07775008964215516
$ws.Range("M6").Value = 44.56267066666667
$ws.Range("N6").Value = 133.688012
$ws.Range("O6").Value = 0.186578549555334
$ws.Range("P6").Value = 0.1926453158454455
$ws.Range("Q6").Value = 269.6577961345924
$ws.Range("R6").Value = 2426.920165211332
$ws.Range("S6").Value = 0.01419496553214443
$ws.Range("T6").Value = 0.01497819057612468

$ws.Range("I7").Value = 0.2215826302097334
$ws.Range("J7").Value = 0.2264456497560767
$ws.Range("M7").Value = 98.170451
$ws.Range("N7").Value = 294.511353
$ws.Range("O7").Value = 0.4110278868558457
$ws.Range("P7").Value = 0.4243928215400083
$ws.Range("Q7").Value = 1730.157500980765
$ws.Range("R7").Value = 15571.41750882688
$ws.Range("S7").Value = 0.09107664025906699
$ws.Range("T7").Value = 0.09610190822544189

$ws.Range("I8").Value = 0.2215826302097334
$ws.Range("J8").Value = 0.2264456497560767
$ws.Range("O8").Value = 0.1453748421312515
$ws.Range("P8").Value = 0.150101833491052
$ws.Range("S8").Value = 0.03221253988576746
$ws.Range("T8").Value = 0.0339899072144597

$ws.Range("I9").Value = 0.2215826302097334
$ws.Range("J9").Value = 0.2264456497560767
$ws.Range("M9").Value = 38.82199566666667
$ws.Range("N9").Value = 116.465987
$ws.Range("O9").Value = 0.1625430328561575
$ws.Range("P9").Value = 0.1678282630971171
$ws.Range("Q9").Value = 684.1994339592683
$ws.Range("R9").Value = 6157.794905633415
$ws.Range("S9").Value = 0.03601671274253448
$ws.Range("T9").Value = 0.03800398008446048

$ws.Range("I10").Value = 0.2215826302097334
$ws.Range("J10").Value = 0.2264456497560767
$ws.Range("M10").Value = 22.5647
$ws.Range("N10").Value = 45.1294
$ws.Range("O10").Value = 0.09447568860141126
$ws.Range("P10").Value = 0.06503176602637677
$ws.Range("Q10").Value = 397.6806112705
$ws.Range("R10").Value = 2386.083667623
$ws.Range("S10").Value = 0.02093417157117643
$ws.Range("T10").Value = 0.01472616051262804

$ws.Range("I11").Value = 0.2215826302097334
$ws.Range("J11").Value = 0.2264456497560767
$ws.Range("M11").Value = 44.56267066666667
$ws.Range("N11").Value = 133.688012
$ws.Range("O11").Value = 0.186578549555334
$ws.Range("P11").Value = 0.1926453158454455
$ws.Range("Q11").Value = 785.3731762693934
$ws.Range("R11").Value = 7068.358586424541
$ws.Range("S11").Value = 0.04134256575118798
$ws.Range("T11").Value = 0.04362369371908653

$ws.Range("G12").Value = 33.62840566666667
$ws.Range("H12").Value = 100.885217
$ws.Range("I12").Value = 0.4228021014155913
$ws.Range("J12").Value = 0.432081235260482
$ws.Range("M12").Value = 98.170451
$ws.Range("N12").Value = 294.511353
$ws.Range("O12").Value = 0.4110278868558457
$ws.Range("P12").Value = 0.4243928215400083
$ws.Range("Q12").Value = 3301.315750707622
$ws.Range("R12").Value = 29711.8417563686
$ws.Range("S12").Value = 0.1737834543030615
$ws.Range("T12").Value = 0.1833721745666881

$ws.Range("G13").Value = 33.62840566666667
$ws.Range("H13").Value = 100.885217
$ws.Range("I13").Value = 0.4228021014155913
$ws.Range("J13").Value = 0.432081235260482
$ws.Range("O13").Value = 0.1453748421312515
$ws.Range("P13").Value = 0.150101833491052
$ws.Range("Q13").Value = 1167.629427180091
$ws.Range("R13").Value = 10508.66484462082
$ws.Range("S13").Value = 0.06146478874605295
$ws.Range("T13").Value = 0.06485618562967693

$ws.Range("G14").Value = 33.62840566666667
$ws.Range("H14").Value = 100.885217
$ws.Range("I14").Value = 0.4228021014155913
$ws.Range("J14").Value = 0.432081235260482
$ws.Range("M14").Value = 38.82199566666667
$ws.Range("N14").Value = 116.465987
$ws.Range("O14").Value = 0.1625430328561575
$ws.Range("P14").Value = 0.1678282630971171
$ws.Range("Q14").Value = 1305.521819068242
$ws.Range("R14").Value = 11749.69637161418
$ws.Range("S14").Value = 0.06872353586204689
$ws.Range("T14").Value = 0.07251544323062353

$ws.Range("G15").Value = 33.62840566666667
$ws.Range("H15").Value = 100.885217
$ws.Range("I15").Value = 0.4228021014155913
$ws.Range("J15").Value = 0.432081235260482
$ws.Range("M15").Value = 22.5647
$ws.Range("N15").Value = 45.1294
$ws.Range("O15").Value = 0.09447568860141126
$ws.Range("P15").Value = 0.06503176602637677
$ws.Range("Q15").Value = 758.8148853466333
$ws.Range("R15").Value = 4552.8893120798
$ws.Range("S15").Value = 0.0399445196733617
$ws.Range("T15").Value = 0.02809900579584752

$ws.Range("G16").Value = 33.62840566666667
$ws.Range("H16").Value = 100.885217
$ws.Range("I16").Value = 0.4228021014155913
$ws.Range("J16").Value = 0.432081235260482
$ws.Range("M16").Value = 44.56267066666667
$ws.Range("N16").Value = 133.688012
$ws.Range("O16").Value = 0.186578549555334
$ws.Range("P16").Value = 0.1926453158454455
$ws.Range("Q16").Value = 1498.571566768734
$ws.Range("R16").Value = 13487.14410091861
$ws.Range("S16").Value = 0.07888580283106823
$ws.Range("T16").Value = 0.08323842603764582

$ws.Range("G17").Value = 5.124275
$ws.Range("H17").Value = 10.24855
$ws.Range("I17").Value = 0.06442631445887793
$ws.Range("J17").Value = 0.04389350863594627
$ws.Range("M17").Value = 98.170451
$ws.Range("N17").Value = 294.511353
$ws.Range("O17").Value = 0.4110278868558457
$ws.Range("P17").Value = 0.4243928215400083
$ws.Range("Q17").Value = 503.052387798025
$ws.Range("R17").Value = 3018.31432678815
$ws.Range("S17").Value = 0.02648101188994282
$ws.Range("T17").Value = 0.01862808997729996

$ws.Range("G18").Value = 5.124275
$ws.Range("H18").Value = 10.24855
$ws.Range("I18").Value = 0.06442631445887793
$ws.Range("J18").Value = 0.04389350863594627
$ws.Range("O18").Value = 0.1453748421312515
$ws.Range("P18").Value = 0.150101833491052
$ws.Range("Q18").Value = 177.92262714655
$ws.Range("R18").Value = 1067.5357628793
$ws.Range("S18").Value = 0.009365965293557743
$ws.Range("T18").Value = 0.006588496124610858

$ws.Range("G19").Value = 5.124275
$ws.Range("H19").Value = 10.24855
$ws.Range("I19").Value = 0.06442631445887793
$ws.Range("J19").Value = 0.04389350863594627
$ws.Range("M19").Value = 38.82199566666667
$ws.Range("N19").Value = 116.465987
$ws.Range("O19").Value = 0.1625430328561575
$ws.Range("P19").Value = 0.1678282630971171
$ws.Range("Q19").Value = 198.9345818448083
$ws.Range("R19").Value = 1193.60749106885
$ws.Range("S19").Value = 0.01047204854789053
$ws.Range("T19").Value = 0.007366571315609172

$ws.Range("G20").Value = 5.124275
$ws.Range("H20").Value = 10.24855
$ws.Range("I20").Value = 0.06442631445887793
$ws.Range("J20").Value = 0.04389350863594627
$ws.Range("M20").Value = 22.5647
$ws.Range("N20").Value = 45.1294
$ws.Range("O20").Value = 0.09447568860141126
$ws.Range("P20").Value = 0.06503176602637677
$ws.Range("Q20").Value = 115.6277280925
$ws.Range("R20").Value = 462.51091237
$ws.Range("S20").Value = 0.006086720422553552
$ws.Range("T20").Value = 0.002854472383689606

$ws.Range("G21").Value = 5.124275
$ws.Range("H21").Value = 10.24855
$ws.Range("I21").Value = 0.06442631445887793
$ws.Range("J21").Value = 0.04389350863594627
$ws.Range("M21").Value = 44.56267066666667
$ws.Range("N21").Value = 133.688012
$ws.Range("O21").Value = 0.186578549555334
$ws.Range("P21").Value = 0.1926453158454455
$ws.Range("Q21").Value = 228.3513792304333
$ws.Range("R21").Value = 1370.1082753826
$ws.Range("S21").Value = 0.01202056830493329
$ws.Range("T21").Value = 0.00845587883473666

$ws.Range("G22").Value = 17.10908866666667
$ws.Range("H22").Value = 51.327266
$ws.Range("I22").Value = 0.2151085815151395
$ws.Range("J22").Value = 0.2198295167053399
$ws.Range("M22").Value = 98.170451
$ws.Range("N22").Value = 294.511353
$ws.Range("O22").Value = 0.4110278868558457
$ws.Range("P22").Value = 0.4243928215400083
$ws.Range("Q22").Value = 1679.606950605656
$ws.Range("R22").Value = 15116.4625554509
$ws.Range("S22").Value = 0.08841562570472622
$ws.Range("T22").Value = 0.0932940688523556

$ws.Range("G23").Value = 17.10908866666667
$ws.Range("H23").Value = 51.327266
$ws.Range("I23").Value = 0.2151085815151395
$ws.Range("J23").Value = 0.2198295167053399
$ws.Range("O23").Value = 0.1453748421312515
$ws.Range("P23").Value = 0.150101833491052
$ws.Range("Q23").Value = 594.0535985396174
$ws.Range("R23").Value = 5346.482386856556
$ws.Range("S23").Value = 0.03127137607884084
$ws.Range("T23").Value = 0.03299681351292336

$ws.Range("G24").Value = 17.10908866666667
$ws.Range("H24").Value = 51.327266
$ws.Range("I24").Value = 0.2151085815151395
$ws.Range("J24").Value = 0.2198295167053399
$ws.Range("M24").Value = 38.82199566666667
$ws.Range("N24").Value = 116.465987
$ws.Range("O24").Value = 0.1625430328561575
$ws.Range("P24").Value = 0.1678282630971171
$ws.Range("Q24").Value = 664.2089660779492
$ws.Range("R24").Value = 5977.880694701542
$ws.Range("S24").Value = 0.03496440123285675
$ws.Range("T24").Value = 0.03689360596613589

$ws.Range("G25").Value = 17.10908866666667
$ws.Range("H25").Value = 51.327266
$ws.Range("I25").Value = 0.2151085815151395
$ws.Range("J25").Value = 0.2198295167053399
$ws.Range("M25").Value = 22.5647
$ws.Range("N25").Value = 45.1294
$ws.Range("O25").Value = 0.09447568860141126
$ws.Range("P25").Value = 0.06503176602637677
$ws.Range("Q25").Value = 386.0614530367334
$ws.Range("R25").Value = 2316.3687182204
$ws.Range("S25").Value = 0.02032253136271561
$ws.Range("T25").Value = 0.01429590169607315

$ws.Range("G26").Value = 17.10908866666667
$ws.Range("H26").Value = 51.327266
$ws.Range("I26").Value = 0.2151085815151395
$ws.Range("J26").Value = 0.2198295167053399
$ws.Range("M26").Value = 44.56267066666667
$ws.Range("N26").Value = 133.688012
$ws.Range("O26").Value = 0.186578549555334
$ws.Range("P26").Value = 0.1926453158454455
$ws.Range("Q26").Value = 762.4266836594659
$ws.Range("R26").Value = 6861.840152935193
$ws.Range("S26").Value = 0.04013464713600005
$ws.Range("T26").Value = 0.04234912667785185
